$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 30 de Marzo de 2020 a las 13:50"

# Update country rows that were reordered / had updated case numbers
$data = @(
    @(18, "Portugal", 6408, 446, 43, 6225, 89, 21, 140),
    @(19, "Canada", 6320, 0, 573, 5682, 120, 0, 65),
    @(20, "Noruega", 4390, 106, 7, 4352, 97, 5, 31),
    @(21, "Israel", 4347, 100, 134, 4197, 80, 1, 16),
    @(22, "Brasil", 4316, 60, 6, 4171, 296, 3, 139),
    @(28, "Dinamarca", 2555, 160, 1, 2477, 113, 5, 77),
    @(30, "Rumania", 1952, 137, 206, 1700, 31, 3, 46),
    @(36, "Pakistan", 1625, 28, 29, 1576, 11, 6, 20),
    @(67, "Lituania", 484, 24, 1, 476, 5, 0, 7),
    @(91, "Vietnam", 203, 9, 52, 151, 3, 0, 0),
    @(94, "Islas Feroe", 168, 9, 70, 98, 1, 0, 0),
    @(95, "Costa de Marfil", 165, 0, 4, 160, 0, 0, 1),
    @(96, "Senegal", 162, 20, 27, 135, 0, 0, 0),
    @(113, "Bolivia", 96, 15, 0, 93, 3, 2, 3),
    @(136, "Uganda", 33, 0, 0, 33, 0, 0, 0),
    @(137, "Barbados", 33, 0, 0, 33, 0, 0, 0),
    @(148, "Tanzania", 19, 5, 1, 18, 0, 0, 0),
    @(149, "Republica de Yibuti", 18, 0, 0, 18, 0, 0, 0),
    @(150, "Islas Virgenes de los Estados Unidos", 17, 0, 0, 17, 0, 0, 0),
    @(151, "Maldivas", 17, 0, 13, 4, 0, 0, 0),
    @(152, "Guinea", 16, 0, 0, 16, 0, 0, 0),
    @(153, "Nueva Caledonia", 15, 0, 0, 15, 0, 0, 0),
    @(154, "Haiti", 15, 0, 1, 14, 0, 0, 0),
    @(156, "Eritrea", 12, 0, 0, 12, 0, 0, 0),
    @(157, "Guinea Ecuatorial", 12, 0, 0, 12, 0, 0, 0),
    @(167, "Santa Lucia", 9, 0, 1, 8, 0, 0, 0),
    @(168, "Siria", 9, 0, 0, 8, 0, 0, 1)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
}

Write-Output "Update complete"
